# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" sheet (copied from the "2022-Q3" sheet's layout)
# right after "总计", trims it down to the single 2022-Q4 holding, and
# refreshes the "总计" summary sheet with the new quarter on top.

function Set-TextValue($ws, $cellRef, $text) {
    # Assigning a numeric-looking string via .Value normally gets coerced to
    # a number (and loses leading zeros, e.g. fund code "004413"). The
    # leading "'" forces Excel to keep it as text; we then restore a plain
    # (unstyled) number format by pasting the format from an untouched
    # neutral cell, so the cell doesn't end up with a stray "@" style.
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range("Z1").Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
}

$wb = $excel.ActiveWorkbook

$zj   = $wb.Worksheets.Item("总计")
$q3   = $wb.Worksheets.Item("2022-Q3")

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet right after "总计", reusing the
#    "2022-Q3" sheet's layout/styles, then overwrite it with the single
#    2022-Q4 fund holding.
# ---------------------------------------------------------------------
$q3.Copy($null, $zj)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

Set-TextValue $q4 "B2" "004413"
$q4.Range("C2").Value = "建信民丰回报定期开放混合"
Set-TextValue $q4 "D2" "0.53"
Set-TextValue $q4 "E2" "20.55"
Set-TextValue $q4 "F2" "0.41"
Set-TextValue $q4 "G2" "0.0022"
$q4.Range("H2").Value = 4

# Only one fund is reported for 2022-Q4 -- drop the extra rows that were
# copied in from 2022-Q3.
$q4.Rows.Item(4).Delete()
$q4.Rows.Item(3).Delete()

# ---------------------------------------------------------------------
# 2. Refresh "总计": push the existing quarters down a row and add the
#    new 2022-Q4 total at the top.
# ---------------------------------------------------------------------
for ($r = 6; $r -ge 2; $r--) {
    $dst = $r + 1
    $zj.Cells.Item($dst, 1).Value = $r - 1
    $zj.Cells.Item($dst, 2).Value = $zj.Cells.Item($r, 2).Value()
    $zj.Cells.Item($dst, 3).Value = $zj.Cells.Item($r, 3).Value()
    $zj.Cells.Item($dst, 4).Value = $zj.Cells.Item($r, 4).Value()
}
$zj.Range("A6").Copy()
$zj.Range("A7").PasteSpecial(-4122)

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q4"
$zj.Range("C2").Value = 1
$zj.Range("D2").Value = 0

# Keep "2021-Q3" (now the last sheet) the selected tab, matching the
# original workbook state.
$q3After = $wb.Worksheets.Item("2021-Q3")
$q3After.Activate()
